# Alligator.Calc.xlsx: rename the "index" column to "i" and switch it from a
# 1-based row counter to a 0-based one (A2 = 0, A3 = 1, ... A503 = 501).
# Also narrow column A now that the header/content is a single character.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "index" -> "i" (this is the header cell of the "testdata" table's
# first column, so the table's column name updates automatically too).
$ws.Range("A1").Value = "i"

# Data rows: 502 rows (A2:A503) change from a 1-based counter (1..502) to a
# 0-based counter (0..501) -- i.e. every value decreases by 1.
for ($r = 2; $r -le 503; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Column A width shrinks from 6 to 4 (character units) to fit the new,
# shorter "i" header/values.
$ws.Columns.Item(1).ColumnWidth = 3.14
